# [FIX] update satellogic consolidated data with interpolation
#
# Inserts three new quarterly rows (Q1-Q3 2022) between the existing
# 2021-12-31 row and the 2022-12-31 row on the "consolidated" sheet,
# linearly interpolating every metric between those two anchor rows,
# and highlights the new quarter-end dates with a yellow fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the old row 3 (2022-12-31 data), which pushes
# that row down to row 6 and leaves rows 3-5 ready for the new quarters.
$ws.Rows("3:5").Insert()

# --- Row 3: 2022-03-31 (25% of the way from row 2 to row 6) ---
$ws.Range("A3").Value = 44681
$ws.Range("B3").Formula = "=(0.25*(B`$6-B`$2))+B`$2"
$ws.Range("C3:S3").Formula = "=(0.25*(C6-C2))+C2"
$ws.Range("U3:AB3").Formula = "=(0.25*(U6-U2))+U2"
$ws.Range("T3").Formula = "=(0.25*(T6-T2))+T2"

# --- Row 4: 2022-06-30 (50% of the way from row 2 to row 6) ---
$ws.Range("A4").Value = 44742
$ws.Range("B4").Formula = "=(0.5*(B`$6-B`$2))+B`$2"
$ws.Range("C4:AB4").Formula = "=(0.5*(C`$6-C`$2))+C`$2"

# --- Row 5: 2022-09-30 (75% of the way from row 2 to row 6) ---
$ws.Range("A5").Value = 44834
$ws.Range("B5").Formula = "=(0.75*(B`$6-B`$2))+B`$2"
$ws.Range("C5:AB5").Formula = "=(0.75*(C`$6-C`$2))+C`$2"

# Highlight the new quarter-end dates (A3:A5) with a yellow fill, matching
# the existing date cell formatting (centered, wrapped, yyyy-mm-dd).
$ws.Range("A3:A5").Interior.Color = 65535

# Move the active selection to A5, as in the edited workbook.
$ws.Range("A5").Select()
